# KARMEN-3: update harmonization rules for TOT_PA_QX, MENOPAUSE and CONTRACEPTIVE
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: TOT_PA_QX ---
# algorithm: (IPAQ_MET_Total*60)/7 -> (IPAQ_MET_Total/60)/7
# status: compatible -> complete
# status_detail: proximate -> compatible
$ws.Range("H6").Value = "(IPAQ_MET_Total/60)/7"
$ws.Range("J6").Value = "complete"
$ws.Range("K6").Value = "compatible"

# --- Row 17: CONTRACEPTIVE ---
# rule_category/algorithm/comment -> impossible; comment cleared
# status -> impossible; status_detail -> incompatible
$ws.Range("F17").Value = "impossible"
$ws.Range("F5").Copy()
$ws.Range("F17").PasteSpecial(-4122)

$ws.Range("G17").Value = "impossible"

$ws.Range("H17").Value = "impossible"

$ws.Range("I17").ClearContents()

$ws.Range("J17").Value = "impossible"
$ws.Range("I5").Copy()
$ws.Range("J17").PasteSpecial(-4122)

$ws.Range("K17").Value = "incompatible"

# --- Row 15: MENOPAUSE ---
# algorithm: recode(0=1;1=2;) -> recode(1=0;2=1;)
$ws.Range("H15").Value = "recode(1=0;2=1;)"

# Update selection to reflect where the author last clicked
$ws.Range("B16").Select()
